$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'328.58"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'6.60%"
$ws.Range("E2").Style = "Normal"
$ws.Range("G2").Value = "'17"
$ws.Range("G2").Style = "Normal"
$ws.Range("D3").Value = "'40.56"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'11.81%"
$ws.Range("E3").Style = "Normal"
$ws.Range("G3").Value = "'17"
$ws.Range("G3").Style = "Normal"
$ws.Range("D4").Value = "'5.900"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'15.40%"
$ws.Range("E4").Style = "Normal"
$ws.Range("G4").Value = "'17"
$ws.Range("G4").Style = "Normal"
$ws.Range("D5").Value = "'0.08131"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'5.32%"
$ws.Range("E5").Style = "Normal"
$ws.Range("G5").Value = "'17"
$ws.Range("G5").Style = "Normal"
$ws.Range("E6").Value = "'4.56%"
$ws.Range("E6").Style = "Normal"
$ws.Range("G6").Value = "'17"
$ws.Range("G6").Style = "Normal"
$ws.Range("E7").Value = "'5.29%"
$ws.Range("E7").Style = "Normal"
$ws.Range("G7").Value = "'17"
$ws.Range("G7").Style = "Normal"
$ws.Range("D8").Value = "'1.953"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'5.85%"
$ws.Range("E8").Style = "Normal"
$ws.Range("G8").Value = "'17"
$ws.Range("G8").Style = "Normal"
$ws.Range("D9").Value = "'2.943"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-0.75%"
$ws.Range("E9").Style = "Normal"
$ws.Range("G9").Value = "'17"
$ws.Range("G9").Style = "Normal"
$ws.Range("D10").Value = "'0.9442"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'2.49%"
$ws.Range("E10").Style = "Normal"
$ws.Range("G10").Value = "'17"
$ws.Range("G10").Style = "Normal"
$ws.Range("D11").Value = "'0.1308"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'15.27%"
$ws.Range("E11").Style = "Normal"
$ws.Range("G11").Value = "'17"
$ws.Range("G11").Style = "Normal"
$ws.Range("D12").Value = "'0.1999"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'8.14%"
$ws.Range("E12").Style = "Normal"
$ws.Range("G12").Value = "'17"
$ws.Range("G12").Style = "Normal"
$ws.Range("D13").Value = "'0.09262"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'6.49%"
$ws.Range("E13").Style = "Normal"
$ws.Range("G13").Value = "'17"
$ws.Range("G13").Style = "Normal"
$ws.Range("D14").Value = "'0.03427"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'3.01%"
$ws.Range("E14").Style = "Normal"
$ws.Range("G14").Value = "'17"
$ws.Range("G14").Style = "Normal"
$ws.Range("D15").Value = "'0.09622"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'1.01%"
$ws.Range("E15").Style = "Normal"
$ws.Range("G15").Value = "'17"
$ws.Range("G15").Style = "Normal"
$ws.Range("D16").Value = "'0.001309"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-5.48%"
$ws.Range("E16").Style = "Normal"
$ws.Range("G16").Value = "'17"
$ws.Range("G16").Style = "Normal"
$ws.Range("D17").Value = "'0.006152"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'1.32%"
$ws.Range("E17").Style = "Normal"
$ws.Range("G17").Value = "'17"
$ws.Range("G17").Style = "Normal"
$ws.Range("B18").Value = 'HotbitToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D18").Value = "'0.004359"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'2.37%"
$ws.Range("E18").Style = "Normal"
$ws.Range("G18").Value = "'17"
$ws.Range("G18").Style = "Normal"
$ws.Range("B19").Value = 'LEO'
$ws.Range("C19").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D19").Value = "'3.372"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-0.03%"
$ws.Range("E19").Style = "Normal"
$ws.Range("G19").Value = "'17"
$ws.Range("G19").Style = "Normal"
$ws.Range("B20").Value = 'BitpandaEcosystemToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D20").Value = "'0.3532"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'2.47%"
$ws.Range("E20").Style = "Normal"
$ws.Range("G20").Value = "'17"
$ws.Range("G20").Style = "Normal"
$ws.Range("B21").Value = 'MCDex'
$ws.Range("C21").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D21").Value = "'7.703"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'21.98%"
$ws.Range("E21").Style = "Normal"
$ws.Range("G21").Value = "'17"
$ws.Range("G21").Style = "Normal"
$ws.Range("B22").Value = 'ProBitToken'
$ws.Range("C22").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D22").Value = "'0.1448"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'11.43%"
$ws.Range("E22").Style = "Normal"
$ws.Range("G22").Value = "'17"
$ws.Range("G22").Style = "Normal"
$ws.Range("B23").Value = 'ZBToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("D23").Value = "'0.2451"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'5.84%"
$ws.Range("E23").Style = "Normal"
$ws.Range("G23").Value = "'17"
$ws.Range("G23").Style = "Normal"
$ws.Range("B24").Value = 'CoinExToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D24").Value = "'0.04447"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'2.41%"
$ws.Range("E24").Style = "Normal"
$ws.Range("G24").Value = "'17"
$ws.Range("G24").Style = "Normal"
$ws.Range("B25").Value = 'BitKan'
$ws.Range("C25").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D25").Value = "'0.001255"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'4.36%"
$ws.Range("E25").Style = "Normal"
$ws.Range("G25").Value = "'17"
$ws.Range("G25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001191"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-10.58%"
$ws.Range("E26").Style = "Normal"
$ws.Range("G26").Value = "'17"
$ws.Range("G26").Style = "Normal"
$ws.Range("D27").Value = "'0.0003994"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'37.47%"
$ws.Range("E27").Style = "Normal"
$ws.Range("G27").Value = "'17"
$ws.Range("G27").Style = "Normal"
$ws.Range("G28").Value = "'17"
$ws.Range("G28").Style = "Normal"
$ws.Range("G29").Value = "'17"
$ws.Range("G29").Style = "Normal"
$ws.Range("G30").Value = "'17"
$ws.Range("G30").Style = "Normal"
$ws.Range("G31").Value = "'17"
$ws.Range("G31").Style = "Normal"
$ws.Range("G32").Value = "'17"
$ws.Range("G32").Style = "Normal"
$ws.Range("G33").Value = "'17"
$ws.Range("G33").Style = "Normal"
$ws.Range("G34").Value = "'17"
$ws.Range("G34").Style = "Normal"
$ws.Range("G35").Value = "'17"
$ws.Range("G35").Style = "Normal"
$ws.Range("G36").Value = "'17"
$ws.Range("G36").Style = "Normal"
$ws.Range("G37").Value = "'17"
$ws.Range("G37").Style = "Normal"
$ws.Range("G38").Value = "'17"
$ws.Range("G38").Style = "Normal"
$ws.Range("D39").Value = "'0.02502"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'19.42%"
$ws.Range("E39").Style = "Normal"
$ws.Range("G39").Value = "'17"
$ws.Range("G39").Style = "Normal"
$ws.Range("D40").Value = "'0.05304"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'8.01%"
$ws.Range("E40").Style = "Normal"
$ws.Range("G40").Value = "'17"
$ws.Range("G40").Style = "Normal"
$ws.Range("D41").Value = "'0.007610"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'0.41%"
$ws.Range("E41").Style = "Normal"
$ws.Range("G41").Value = "'17"
$ws.Range("G41").Style = "Normal"
$ws.Range("D42").Value = "'0.1431"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'6.18%"
$ws.Range("E42").Style = "Normal"
$ws.Range("G42").Value = "'17"
$ws.Range("G42").Style = "Normal"
$ws.Range("D43").Value = "'0.008974"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'5.66%"
$ws.Range("E43").Style = "Normal"
$ws.Range("G43").Value = "'17"
$ws.Range("G43").Style = "Normal"
$ws.Range("D44").Value = "'0.002069"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-0.30%"
$ws.Range("E44").Style = "Normal"
$ws.Range("G44").Value = "'17"
$ws.Range("G44").Style = "Normal"
$ws.Range("D45").Value = "'0.009494"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'12.44%"
$ws.Range("E45").Style = "Normal"
$ws.Range("G45").Value = "'17"
$ws.Range("G45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006856"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'6.18%"
$ws.Range("E46").Style = "Normal"
$ws.Range("G46").Value = "'17"
$ws.Range("G46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-0.07%"
$ws.Range("E47").Style = "Normal"
$ws.Range("G47").Value = "'17"
$ws.Range("G47").Style = "Normal"
$ws.Range("D48").Value = "'0.002899"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-12.18%"
$ws.Range("E48").Style = "Normal"
$ws.Range("G48").Value = "'17"
$ws.Range("G48").Style = "Normal"
$ws.Range("D49").Value = "'0.001802"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'24.65%"
$ws.Range("E49").Style = "Normal"
$ws.Range("G49").Value = "'17"
$ws.Range("G49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.07%"
$ws.Range("E50").Style = "Normal"
$ws.Range("G50").Value = "'17"
$ws.Range("G50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'-0.07%"
$ws.Range("E51").Style = "Normal"
$ws.Range("G51").Value = "'17"
$ws.Range("G51").Style = "Normal"
